# HP added to weekly content
#
# The "3. if/elif/else" exercise bullet (cell D3) is expanded to mention the
# new "earth layers" exercise and a dedicated "nested conditionals" line, and
# the sheet's active selection moves from C2 to D4 (with row 4 scrolled into
# view) to reflect where the author was last working.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "if/elif/else" exercises cell (D3) with the new wording.
$newD3 = "1. A set of if statements (point to where to find more help in the hint)`n" + `
         "2. A set of if/else statements`n" + `
         "3. if/elif/else (earth layers)`n" + `
         "nested conditionals `n" + `
         "4. multiple elif (what's the correct order?)`n" + `
         "5. Open ended question - text based adventure game - get a friend or TA to play it or send it to me! "

$ws.Range("D3").Value = $newD3

# Move the selection to D4 (and scroll so row 4 is visible), matching the
# author's cursor position when the workbook was saved.
$ws.Range("D4").Select()

Write-Host "Updated D3 wording and moved selection to D4"
